$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 SmartScore cells from text to numeric ---
$ws.Range('I5').Value2 = 0.6
$ws.Range('L5').Value2 = 0.563
$ws.Range('O5').Value2 = 0.539
$ws.Range('R5').Value2 = 0.643
$ws.Range('U5').Value2 = 0.582
$ws.Range('X5').Value2 = 0.569
$ws.Range('AA5').Value2 = 0.624
$ws.Range('AD5').Value2 = 0.611
$ws.Range('AG5').Value2 = 0.611

# --- Add new row 6 with participant data (Hanna Moriel) ---

# D6 numeric age
$ws.Range("D6").Value2 = 21

# Plain text cells
$ws.Range('A6').Value2 = 'Hanna Moriel_20251113_214222'
$ws.Range('C6').Value2 = 'Hanna Moriel'
$ws.Range('E6').Value2 = 'Female'
$ws.Range('F6').Value2 = '2025-11-13 21:42:22'
$ws.Range('H6').Value2 = 'Maruchan Ramen Sabor Pollo'
$ws.Range('J6').Value2 = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range('K6').Value2 = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range('M6').Value2 = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range('N6').Value2 = 'Nongshim Neoguri Spicy Seafood'
$ws.Range('P6').Value2 = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range('Q6').Value2 = 'Kraft Macaroni & Cheese Dinner'
$ws.Range('S6').Value2 = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range('T6').Value2 = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range('V6').Value2 = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range('W6').Value2 = 'Velveeta Original Shells & Cheese (microwave cups)'
$ws.Range('Y6').Value2 = 'Muy cremoso, porción individual, rápido, salado, ideal para niños'
$ws.Range('Z6').Value2 = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range('AB6').Value2 = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range('AC6').Value2 = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range('AE6').Value2 = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range('AF6').Value2 = 'Jack Link’s Beef Jerky Original'
$ws.Range('AH6').Value2 = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'

# B6: empty inline string cell (kept blank, matching source)
$ws.Range("B6").Value2 = ""

# G6: JSON weights text (multi-line)
$json6 = @"
{
  "portion": 0.4,
  "diet": 0.2857142857142857,
  "salt": 0.8,
  "fat": 0.4,
  "natural": 0.2,
  "convenience": 0.4,
  "price": 0.2
}
"@
$ws.Range("G6").Value2 = $json6

# SmartScore cells in row 6 remain TEXT (e.g. "0.540"), force text format so Excel
# does not auto-convert the numeric-looking string to a number, then restore the
# default "Normal" style so no stray per-cell style is left behind.
$ws.Range('I6').NumberFormat = "@"
$ws.Range('I6').Value2 = '0.540'
$ws.Range('I6').Style = "Normal"
$ws.Range('L6').NumberFormat = "@"
$ws.Range('L6').Value2 = '0.520'
$ws.Range('L6').Style = "Normal"
$ws.Range('O6').NumberFormat = "@"
$ws.Range('O6').Value2 = '0.449'
$ws.Range('O6').Style = "Normal"
$ws.Range('R6').NumberFormat = "@"
$ws.Range('R6').Value2 = '0.622'
$ws.Range('R6').Style = "Normal"
$ws.Range('U6').NumberFormat = "@"
$ws.Range('U6').Value2 = '0.616'
$ws.Range('U6').Style = "Normal"
$ws.Range('X6').NumberFormat = "@"
$ws.Range('X6').Value2 = '0.615'
$ws.Range('X6').Style = "Normal"
$ws.Range('AA6').NumberFormat = "@"
$ws.Range('AA6').Value2 = '0.729'
$ws.Range('AA6').Style = "Normal"
$ws.Range('AD6').NumberFormat = "@"
$ws.Range('AD6').Value2 = '0.702'
$ws.Range('AD6').Style = "Normal"
$ws.Range('AG6').NumberFormat = "@"
$ws.Range('AG6').Value2 = '0.685'
$ws.Range('AG6').Style = "Normal"

Write-Host "Done updating sheet."
